# Atualizando o andamento das SARs.
#
# Applies to the SAR03 table:
#   - "Data da ultima alteracao": 20/05/2015 -> 21/05/2015 (and moves the
#     Word "_GoBack" last-edit bookmark here, as real Word would after the
#     most recent keystrokes)
#   - "Descricao": removes the now-stale "_GoBack" bookmark left over from a
#     previous edit session (text itself is unchanged)
#   - "Impacto": re-unifies the paragraph's runs (text unchanged)
#   - "Status da Alteracao": "Em avaliação." -> "Aprovada para resolução."

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# 1) "Data da ultima alteracao" (row 4, col 2): 20/05/2015 -> 21/05/2015
# ---------------------------------------------------------------------
$dateCell = $t.Cell(4, 2)
$dateStart = $dateCell.Range.Start

$d.Range($dateStart, $dateStart + 1).Text = "2"
$d.Range($dateStart + 1, $dateStart + 2).Text = "1"

# ---------------------------------------------------------------------
# 2) "Descricao" (row 5, col 2): text is already correct, but it still
#    carries the old "_GoBack" bookmark splitting the runs; force a real
#    replace (via a throwaway placeholder) so the bookmark is cleared and
#    the paragraph collapses back to a single run.
# ---------------------------------------------------------------------
$descCell = $t.Cell(5, 2)
$descRange = $d.Range($descCell.Range.Start, $descCell.Range.End - 1)
$descFinal = "Remoção do atributo “_id da Conta” da entidade “atividade”."
$descRange.Text = "@@TMP@@"
$d.Range($descCell.Range.Start, $descCell.Range.Start + 7).Text = $descFinal

# ---------------------------------------------------------------------
# 3) "Impacto" (row 7, col 2): text unchanged, merge the 3 runs back into
#    one the same way.
# ---------------------------------------------------------------------
$impactoCell = $t.Cell(7, 2)
$impactoRange = $d.Range($impactoCell.Range.Start, $impactoCell.Range.End - 1)
$impactoFinal = "Remover este atributo eliminaria uma redundância de informação evidente. Não há impacto negativo para esta alteração."
$impactoRange.Text = "@@TMP@@"
$d.Range($impactoCell.Range.Start, $impactoCell.Range.Start + 7).Text = $impactoFinal

# ---------------------------------------------------------------------
# 4) "Status da Alteracao" (row 9, col 2): Em avaliação. -> Aprovada para
#    resolução.
# ---------------------------------------------------------------------
$statusCell = $t.Cell(9, 2)
$statusRange = $d.Range($statusCell.Range.Start, $statusCell.Range.End - 1)
$statusRange.Text = "Aprovada para resolução."

# ---------------------------------------------------------------------
# 5) Move the single "_GoBack" bookmark: drop the stale one (in
#    Descricao - already gone, cleared by step 2) and add the fresh one
#    at the site of the most recent edit (between "21" and "/05/2015").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$gbPos = $dateStart + 2
$d.Bookmarks.Add("_GoBack", $d.Range($gbPos, $gbPos))

Write-Output "SAR03 status update applied."
